$d = $word.ActiveDocument

# --- Change 1: "Integrar el sistema de pago mediante la plataforma TRANSBANK." -> split into runs with "Khipu" ---
$r1 = $d.Content
$ok1 = $r1.Find.Execute("Integrar el sistema de pago mediante la plataforma TRANSBANK.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok1) { throw "Change 1 target text not found" }
$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="7A6AE38E" w14:textId="77777777" w:rsidR="009C36C0" w:rsidRPr="002F1E81" w:rsidRDefault="009C36C0" w:rsidP="009C36C0"><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="8"/></w:numPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:i/><w:sz w:val="18"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-CL"/></w:rPr></w:pPr><w:r w:rsidRPr="002F1E81"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:i/><w:sz w:val="18"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-CL"/></w:rPr><w:t>Integrar el sistema de pago mediante la plataforma</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:i/><w:sz w:val="18"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-CL"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:i/><w:sz w:val="18"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-CL"/></w:rPr><w:t>Khipu</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:i/><w:sz w:val="18"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-CL"/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r1.InsertXML($xml1)

# --- Change 2: "Integración de pagos: Se desarrolló ... TRANSBANK para facilitar ..." -> split into runs with "Khipu" ---
$r2 = $d.Content
$ok2 = $r2.Find.Execute(": Se desarrolló la conexión con la plataforma TRANSBANK para facilitar pagos electrónicos, permitiendo a los usuarios realizar transacciones de manera rápida y segura a través de la aplicación.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok2) { throw "Change 2 target text not found" }
$r2full = $d.Range($r2.Paragraphs(1).Range.Start, $r2.Paragraphs(1).Range.End)
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="40C2FBDA" w14:textId="77777777" w:rsidR="009C36C0" w:rsidRPr="002F1E81" w:rsidRDefault="009C36C0" w:rsidP="009C36C0"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:i/><w:sz w:val="18"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-CL"/></w:rPr></w:pPr><w:r w:rsidRPr="002F1E81"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:b/><w:bCs/><w:i/><w:sz w:val="18"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-CL"/></w:rPr><w:lastRenderedPageBreak/><w:t>Integración de pagos</w:t></w:r><w:r w:rsidRPr="002F1E81"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:i/><w:sz w:val="18"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-CL"/></w:rPr><w:t xml:space="preserve">: Se desarrolló la conexión con la plataforma </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:i/><w:sz w:val="18"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-CL"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:i/><w:sz w:val="18"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-CL"/></w:rPr><w:t>Khipu</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:i/><w:sz w:val="18"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-CL"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:i/><w:sz w:val="18"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-CL"/></w:rPr><w:t>para facilitar pagos electrónicos, permitiendo a los usuarios realizar transacciones de manera rápida y segura a través de la aplicación.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r2full.InsertXML($xml2)

# --- Change 3: merge "Mockups ... app ... y de cómo ..." runs, removing gramStart/gramEnd proofErr wrapping ---
$r3 = $d.Content
$ok3 = $r3.Find.Execute("Las evidencias que tenemos son los Mockups de la aplicación de cómo se ha llevado a cabo la ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok3) { throw "Change 3 target text not found" }
$r3full = $d.Range($r3.Paragraphs(1).Range.Start, $r3.Paragraphs(1).Range.End)
$xml3 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="4B8D460D" w14:textId="72AD552F" w:rsidR="002F1E81" w:rsidRPr="003605F9" w:rsidRDefault="002F1E81" w:rsidP="002F1E81"><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="12"/></w:numPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:i/><w:sz w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="003605F9"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:i/><w:sz w:val="18"/></w:rPr><w:t xml:space="preserve">Evidencia: </w:t></w:r><w:r w:rsidRPr="003605F9"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:i/><w:sz w:val="18"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-CL"/></w:rPr><w:t xml:space="preserve">Las evidencias que tenemos son los Mockups de la aplicación de cómo se ha llevado a cabo la app y de cómo hemos ido avanzando con el proyecto con todas sus funciones correctas sin ningún problema o error, </w:t></w:r><w:r w:rsidRPr="003605F9"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:i/><w:sz w:val="18"/></w:rPr><w:t>que muestran las primeras funcionalidades administrativas del sistema en desarrollo, tales como la interfaz de usuario, panel de control, y algunas interacciones básicas.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r3full.InsertXML($xml3)

Write-Output "done"
